$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.393.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3760"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3422"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.016"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.952"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.580.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06734"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.223"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.399.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.394"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.681"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.08%  "
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.032"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.749.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.153"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9861"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.973"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08481"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02542"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.388"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2314"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06575"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.418"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6400"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.789"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5971"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.300"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.090"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("E51").Value = "  +0.64%  "
